# Fix registration and login
# - rotate the stored password hash for the existing user (row 2)
# - add a new registered user (row 3): M'mah Zombo

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the existing user's password hash (D2)
$ws.Range("D2").Value = '$2a$10$Bz4/5bkmPGFx.KNDLW2Us.iO2Q9dDQjt0wGkqqVFHyeLfHTjb.EF.'

# Append the new user row
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "M'mah Zombo"
$ws.Range("C3").Value = "zombo@gmail.com"
$ws.Range("D3").Value = '$2a$10$/sNF3JT9o2N3GlJj//AFE.bxtu9fT9CyTXrZD1iaVRM9g9nH8UaCa'
$ws.Range("E3").Value = "agent"

# avatar column (F) is blank for every user; "'" types an empty, explicit
# text value (like Excel's leading apostrophe) instead of clearing the cell
$ws.Range("F3").Value = "'"
$ws.Range("F3").Style = "Normal"
